$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename header "No Pengumuman" -> "Tanggal Pendataan" (column D)
$ws.Range("D1").Value = "Tanggal Pendataan"

# Fill in the previously empty "Tanggal Pendataan" value for the data row.
# Force text interpretation so the dd-mm-yyyy-looking string isn't parsed
# into a date serial number, then clear the formatting override so the
# cell keeps the workbook's default (unstyled) look.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "11-02-2023"
$ws.Range("D2").ClearFormats()

# Reformat "Tanggal Lahir" / "An Tanggal Lahir" from yyyy-mm-dd to dd-mm-yyyy
$ws.Range("O2").Value = "20-08-2000"
$ws.Range("W2").Value = "20-08-2000"

# Fix "Desa Letak Tanah" spelling: PONDOKJOYO -> PONDOK JOYO
$ws.Range("AC2").Value = "PONDOK JOYO"
